$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Month_Start"
$ws.Range("C1").Value = "Month_End"

$ws.Range("D3").Select()
